$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of mod-count data (row 27).
# Column A must stay a literal text string like "2025/12/06" rather than
# being auto-converted into a date serial number, so we temporarily force
# a text number format while assigning it, then restore the default
# ("Normal") style before re-applying the same centered alignment used by
# the rest of the data rows (e.g. row 26).
$ws.Cells.Item(27, 1).NumberFormat = "@"
$ws.Cells.Item(27, 1).Value = "2025/12/06"
$ws.Cells.Item(27, 1).Style = "Normal"

$ws.Cells.Item(27, 2).Value = "逃离鸭科夫"
$ws.Cells.Item(27, 3).Value = 1336

# Match formatting/style of the preceding data row (row 26), which is
# centered both horizontally and vertically.
$ws.Range("A27:C27").HorizontalAlignment = -4108
$ws.Range("A27:C27").VerticalAlignment = -4108
